$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3's value from the stale code to the corrected one
$ws.Range("B3").Value = "zaz"

# Remove the now-obsolete duplicate rows (4-7) that held old fetch data
$ws.Range("A4:B7").ClearContents()
